$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Henderson")
$ws.Activate()

# Fill in missing age-sample values for 2019 (row 57, Fence sample)
$ws.Range("E57").Value = 16
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 1
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0

# Fill in missing age-sample size for 2020 (row 58, Fence sample)
$ws.Range("E58").Value = 0

# Insert a new row before row 61 (old "2023/Catch" row) to add a
# new "2022/Fence sample" record, pushing the old row 61 down to 62
$ws.Rows("61").Insert()

$ws.Range("A61").Value = 2022
$ws.Range("B61").Value = 18646
$ws.Range("C61").Value = 7731
$ws.Range("D61").Value = "Fence sample"
$ws.Range("E61").Value = 15
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0.33333333333333331
$ws.Range("H61").Value = 0.66666666666666663
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0

# Add new row 63 for 2023 "Fence sample" record
$ws.Range("A63").Value = 2023
$ws.Range("B63").Value = 13113
$ws.Range("C63").Value = 8196
$ws.Range("D63").Value = "Fence sample"
$ws.Range("E63").Value = 7
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0.14285714285714285
$ws.Range("H63").Value = 0.8571428571428571
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0

# Scroll so row 34 is at the top of the view and select I64, matching
# the author's final on-screen position after the edits
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("I64").Select()
